# Auto-generated: update Sheets via scheduled runner
# Refreshes market-price derived columns (H:N) across several sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 163.42857
$ws.Range("I5").Value = 163.42857
$ws.Range("K5").Value = 163.42857
$ws.Range("M5").Value = -48.42857000000001
$ws.Range("H32").Value = 5871.7144
$ws.Range("I32").Value = 6000.3335
$ws.Range("J32").Value = 5775.25
$ws.Range("K32").Value = 6000.3335
$ws.Range("L32").Value = 5775.25
$ws.Range("M32").Value = -5674.3335
$ws.Range("N32").Value = -6427.25
$ws.Range("H76").Value = 4006.6365
$ws.Range("I76").Value = 3407.3
$ws.Range("J76").Value = 10000
$ws.Range("K76").Value = 3407.3
$ws.Range("L76").Value = 10000
$ws.Range("M76").Value = -3092.3
$ws.Range("N76").Value = -10630
$ws.Range("H79").Value = 4006.6365
$ws.Range("I79").Value = 3407.3
$ws.Range("J79").Value = 10000
$ws.Range("K79").Value = 3407.3
$ws.Range("L79").Value = 10000
$ws.Range("M79").Value = -2315.3
$ws.Range("N79").Value = -12184
$ws.Range("H88").Value = 1673.5294
$ws.Range("I88").Value = 1799.4445
$ws.Range("J88").Value = 1531.875
$ws.Range("K88").Value = 1799.4445
$ws.Range("L88").Value = 1531.875
$ws.Range("M88").Value = -1393.4445
$ws.Range("N88").Value = -2343.875
$ws.Range("H91").Value = 1673.5294
$ws.Range("I91").Value = 1799.4445
$ws.Range("J91").Value = 1531.875
$ws.Range("K91").Value = 1799.4445
$ws.Range("L91").Value = 1531.875
$ws.Range("M91").Value = -395.4445000000001
$ws.Range("N91").Value = -4339.875
$ws.Range("H98").Value = 543.8946999999999
$ws.Range("I98").Value = 546.3333
$ws.Range("K98").Value = 546.3333
$ws.Range("M98").Value = 951.6667
$ws.Range("H122").Value = 543.8946999999999
$ws.Range("I122").Value = 546.3333
$ws.Range("K122").Value = 1638.9999
$ws.Range("M122").Value = 811.0001
$ws.Range("H131").Value = 1400
$ws.Range("I131").Value = 1400
$ws.Range("K131").Value = 4200
$ws.Range("M131").Value = 840
$ws.Range("H135").Value = 41202.44
$ws.Range("I135").Value = 1097.4445
$ws.Range("K135").Value = 9877.0005
$ws.Range("M135").Value = -7342.0005
$ws.Range("H137").Value = 3504.568
$ws.Range("I137").Value = 1324.2727
$ws.Range("J137").Value = 10045.454
$ws.Range("K137").Value = 3972.8181
$ws.Range("L137").Value = 30136.362
$ws.Range("M137").Value = -1422.8181
$ws.Range("N137").Value = -35236.362
$ws.Range("H141").Value = 52354.277
$ws.Range("I141").Value = 52354.277
$ws.Range("K141").Value = 157062.831
$ws.Range("M141").Value = -151882.831

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 176.875
$ws.Range("I5").Value = 143.4
$ws.Range("K5").Value = 143.4
$ws.Range("M5").Value = -31.40000000000001
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H61").Value = 1274.4584
$ws.Range("I61").Value = 908.6667
$ws.Range("J61").Value = 2371.8333
$ws.Range("K61").Value = 908.6667
$ws.Range("L61").Value = 2371.8333
$ws.Range("M61").Value = -696.6667
$ws.Range("N61").Value = -2795.8333
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 12987.556
$ws.Range("I132").Value = 15519.643
$ws.Range("J132").Value = 4125.25
$ws.Range("K132").Value = 46558.929
$ws.Range("L132").Value = 12375.75
$ws.Range("M132").Value = -44028.929
$ws.Range("N132").Value = -17435.75
$ws.Range("H136").Value = 1274.4584
$ws.Range("I136").Value = 908.6667
$ws.Range("J136").Value = 2371.8333
$ws.Range("K136").Value = 2726.0001
$ws.Range("L136").Value = 7115.499899999999
$ws.Range("M136").Value = -176.0001000000002
$ws.Range("N136").Value = -12215.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 176.875
$ws.Range("I4").Value = 143.4
$ws.Range("K4").Value = 143.4
$ws.Range("M4").Value = -28.40000000000001
$ws.Range("H20").Value = 8935.875
$ws.Range("I20").Value = 6958.8887
$ws.Range("J20").Value = 11477.714
$ws.Range("K20").Value = 6958.8887
$ws.Range("L20").Value = 11477.714
$ws.Range("M20").Value = -6711.8887
$ws.Range("N20").Value = -11971.714
$ws.Range("H26").Value = 14839.667
$ws.Range("I26").Value = 7000
$ws.Range("K26").Value = 7000
$ws.Range("M26").Value = -6708
$ws.Range("H106").Value = 32000
$ws.Range("J106").Value = 32000
$ws.Range("L106").Value = 32000
$ws.Range("N106").Value = -34524
$ws.Range("H134").Value = 1370.64
$ws.Range("I134").Value = 1376.174
$ws.Range("K134").Value = 4128.522
$ws.Range("M134").Value = -1593.522

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 46896.74
$ws.Range("J122").Value = 3630.1428
$ws.Range("L122").Value = 10890.4284
$ws.Range("N122").Value = -15790.4284
$ws.Range("H134").Value = 5762.2
$ws.Range("I134").Value = 5952.75
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 17858.25
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -15323.25
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 111156550
$ws.Range("J37").Value = 111156550
$ws.Range("L37").Value = 333469650
$ws.Range("N37").Value = -333469874
$ws.Range("H116").Value = 1878.5
$ws.Range("I116").Value = 1410
$ws.Range("J116").Value = 7032
$ws.Range("K116").Value = 4230
$ws.Range("L116").Value = 21096
$ws.Range("M116").Value = -788
$ws.Range("N116").Value = -27980
$ws.Range("H118").Value = 2409.5454
$ws.Range("I118").Value = 250.83333
$ws.Range("K118").Value = 752.49999
$ws.Range("M118").Value = 490.50001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H70").Value = 5672
$ws.Range("I70").Value = 5008
$ws.Range("J70").Value = 7000
$ws.Range("K70").Value = 5008
$ws.Range("L70").Value = 7000
$ws.Range("M70").Value = -4738
$ws.Range("N70").Value = -7540
$ws.Range("H73").Value = 5672
$ws.Range("I73").Value = 5008
$ws.Range("J73").Value = 7000
$ws.Range("K73").Value = 5008
$ws.Range("L73").Value = 7000
$ws.Range("M73").Value = -4072
$ws.Range("N73").Value = -8872
$ws.Range("H122").Value = 314653.3
$ws.Range("I122").Value = 359032.44
$ws.Range("K122").Value = 1077097.32
$ws.Range("M122").Value = -1074647.32
$ws.Range("H126").Value = 2906.889
$ws.Range("I126").Value = 2341.3333
$ws.Range("J126").Value = 4038
$ws.Range("K126").Value = 7023.999899999999
$ws.Range("L126").Value = 12114
$ws.Range("M126").Value = -4553.999899999999
$ws.Range("N126").Value = -17054
$ws.Range("H132").Value = 3922.5
$ws.Range("I132").Value = 3776.5
$ws.Range("J132").Value = 4141.5
$ws.Range("K132").Value = 11329.5
$ws.Range("L132").Value = 12424.5
$ws.Range("M132").Value = -8799.5
$ws.Range("N132").Value = -17484.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3329.625
$ws.Range("I7").Value = 2326.6
$ws.Range("K7").Value = 2326.6
$ws.Range("M7").Value = -2214.6
$ws.Range("H22").Value = 5000
$ws.Range("J22").Value = 5000
$ws.Range("L22").Value = 5000
$ws.Range("N22").Value = -5590
$ws.Range("H27").Value = 5000
$ws.Range("J27").Value = 5000
$ws.Range("L27").Value = 5000
$ws.Range("N27").Value = -5214
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H46").Value = 29940.934
$ws.Range("I46").Value = 34162.617
$ws.Range("J46").Value = 2500
$ws.Range("K46").Value = 34162.617
$ws.Range("L46").Value = 2500
$ws.Range("M46").Value = -33974.617
$ws.Range("N46").Value = -2876
$ws.Range("H82").Value = 1234.7858
$ws.Range("J82").Value = 1098.5555
$ws.Range("L82").Value = 1098.5555
$ws.Range("N82").Value = -1820.5555
$ws.Range("H85").Value = 1234.7858
$ws.Range("J85").Value = 1098.5555
$ws.Range("L85").Value = 1098.5555
$ws.Range("N85").Value = -3594.5555
$ws.Range("H122").Value = 4045.8572
$ws.Range("I122").Value = 2955.25
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 8865.75
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -6415.75
$ws.Range("N122").Value = -21400
$ws.Range("H126").Value = 3329.625
$ws.Range("I126").Value = 2326.6
$ws.Range("K126").Value = 6979.799999999999
$ws.Range("M126").Value = -4509.799999999999
$ws.Range("H135").Value = 10000
$ws.Range("I135").Value = 10000
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 10000
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -4930
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 2703.92
$ws.Range("I136").Value = 2304.682
$ws.Range("J136").Value = 5631.6665
$ws.Range("K136").Value = 6914.045999999999
$ws.Range("L136").Value = 16894.9995
$ws.Range("M136").Value = -4364.045999999999
$ws.Range("N136").Value = -21994.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3229.762
$ws.Range("J81").Value = 1083.1666
$ws.Range("L81").Value = 2166.3332
$ws.Range("N81").Value = -4288.3332
$ws.Range("H84").Value = 3229.762
$ws.Range("J84").Value = 1083.1666
$ws.Range("L84").Value = 10831.666
$ws.Range("N84").Value = -21439.666
$ws.Range("H126").Value = 2970
$ws.Range("I126").Value = 1946.5
$ws.Range("K126").Value = 5839.5
$ws.Range("M126").Value = -3369.5
$ws.Range("H132").Value = 17102.11
$ws.Range("I132").Value = 20652.875
$ws.Range("K132").Value = 61958.625
$ws.Range("M132").Value = -59428.625
$ws.Range("H136").Value = 2701.2856
$ws.Range("I136").Value = 2786.5186
$ws.Range("J136").Value = 400
$ws.Range("K136").Value = 8359.5558
$ws.Range("L136").Value = 1200
$ws.Range("M136").Value = -5809.5558
$ws.Range("N136").Value = -6300
